# Use Case Login Logout Create course sequence.docx
# Applies the edits described by the commit:
# "#39 Login Logout Create Crouse sequence requirements.
#  Those use case were done by Marc-Andre Leclair"

$d = $word.ActiveDocument

function Get-CellRangeNoMark($cell) {
    $r = $cell.Range
    $r.End = $r.End - 1
    return $r
}

# ---------------------------------------------------------------------
# Table 1 : "Login"
# ---------------------------------------------------------------------
$tbl1 = $d.Tables(1)

# Minimum Guarantee row (9) -> append " User will not log in"
$cell = $tbl1.Cell(9, 1)
$r = Get-CellRangeNoMark $cell
$r.Text = "Minimum Guarantee:   User will not log in"

# Importance assessment row (12) -> append " 5"
$cell = $tbl1.Cell(12, 1)
$r = Get-CellRangeNoMark $cell
$r.Text = "Importance assessment: 5"

# ---------------------------------------------------------------------
# Insert a new summary paragraph right after Table 1, before "Case: "
# ---------------------------------------------------------------------
$anchor = $d.Range($tbl1.Range.End, $tbl1.Range.End).Paragraphs(1)
$anchor.Range.InsertParagraphAfter()
$newPara = $anchor.Next()
$npr = $newPara.Range
$npr.End = $npr.End - 1
$npr.Text = " This Use Case enables any user to log on onto the scheduler to access their own profile to then effectuate whatever they need to in order to accomplish what they initially needed to."

# ---------------------------------------------------------------------
# Table 2 : "Logout"
# ---------------------------------------------------------------------
$tbl2 = $d.Tables(2)

# Minimum Guarantee row (9) -> append " User will log out"
$cell = $tbl2.Cell(9, 1)
$r = Get-CellRangeNoMark $cell
$r.Text = "Minimum Guarantee:   User will log out"

# Importance assessment row (12) -> append "5" (no leading space)
$cell = $tbl2.Cell(12, 1)
$r = Get-CellRangeNoMark $cell
$r.Text = "Importance assessment:5"

# Traces to test Case row (13) -> append " "
$cell = $tbl2.Cell(13, 1)
$r = Get-CellRangeNoMark $cell
$r.Text = "Traces to test Case: "

# ---------------------------------------------------------------------
# Replace the 2nd of the 6 empty paragraphs after Table 2 with the
# TimeTurner paragraph (includes proofErr spell-check markers)
# ---------------------------------------------------------------------
$p1 = $d.Range($tbl2.Range.End, $tbl2.Range.End).Paragraphs(1)
$p2 = $p1.Next()
$p2r = $p2.Range
$p2r.End = $p2r.End - 1
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' +
        '<w:r><w:t xml:space="preserve">   When the user is done his work on his </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>TimeTurner</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> profile, he or she will be able to log out from their account, cutting all their connection with any data associated to their username.</w:t></w:r>' +
        '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2r.InsertXML($xml2)

# ---------------------------------------------------------------------
# Table 3 : "Create course sequence"
# ---------------------------------------------------------------------
$tbl3 = $d.Tables(3)

# Use Case Name row (1) -> merge runs, no text change
$cell = $tbl3.Cell(1, 1)
$r = Get-CellRangeNoMark $cell
$r.Text = "Use Case Name:       Create course sequence"

# Actor(s) row (4) -> merge runs, no text change
$cell = $tbl3.Cell(4, 1)
$r = Get-CellRangeNoMark $cell
$r.Text = "Actor(s):                     Primary Actor(s): Student"

# Goal/Actor Goals row (5) -> merge runs, no text change
$cell = $tbl3.Cell(5, 1)
$r = Get-CellRangeNoMark $cell
$r.Text = "Goal/ Actor Goals:   A student wants to create a personalize course sequence"

# Description/Summary row (6) -> merge runs, no text change
$cell = $tbl3.Cell(6, 1)
$r = Get-CellRangeNoMark $cell
$r.Text = "Description/Summary:   The Student wants to personalize his course sequence by adding/ removing courses from its sequence to personalize it to his or her own need."

# Preconditions row (7) -> merge runs in first paragraph only, no text change
$cell = $tbl3.Cell(7, 1)
$firstPara = $cell.Range.Paragraphs(1)
$fpr = $firstPara.Range
$fpr.End = $fpr.End - 1
$fpr.Text = "Preconditions:    The user is authenticated."

# Post-conditions row (8) -> merge runs in first paragraph only, no text change
$cell = $tbl3.Cell(8, 1)
$firstPara = $cell.Range.Paragraphs(1)
$fpr = $firstPara.Range
$fpr.End = $fpr.End - 1
$fpr.Text = "Post-conditions:   The user is now enrolled / has a personalize course sequence for the future"

# Minimum Guarantee row (9) -> append " No course sequence will be created"
$cell = $tbl3.Cell(9, 1)
$r = Get-CellRangeNoMark $cell
$r.Text = "Minimum Guarantee:   No course sequence will be created"

# Basic Flow row (10) -> merge runs in first paragraph, no text change;
# also remove the _GoBack bookmark from the last paragraph of this cell
$cell = $tbl3.Cell(10, 1)
$firstPara = $cell.Range.Paragraphs(1)
$fpr = $firstPara.Range
$fpr.End = $fpr.End - 1
$fpr.Text = "Basic Flow:                1.Student request to make a change to its sequence"

$lastPara = $cell.Range.Paragraphs(6)
$lpr = $lastPara.Range
$lpr.End = $lpr.End - 1
$lpr.Text = "                                    6. The course sequence is updated and shown to the student."

# Importance assessment row (12) -> append " 5"
$cell = $tbl3.Cell(12, 1)
$r = Get-CellRangeNoMark $cell
$r.Text = "Importance assessment: 5"

# ---------------------------------------------------------------------
# Replace the single trailing empty paragraph after Table 3 with two new
# paragraphs: a single space, then the long closing paragraph that
# contains the (moved) _GoBack bookmark and the TimeTurner proofErr span.
# ---------------------------------------------------------------------
$tailPara = $d.Range($tbl3.Range.End, $tbl3.Range.End).Paragraphs(1)
$tpr = $tailPara.Range
$tpr.End = $tpr.End - 1
$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' +
        '<w:p>' +
        '<w:r><w:t xml:space="preserve">Once a student is logged in, they can create a course sequence. This course sequence will give the student preferences upon when they want their classes. For instance, this could include morning versus evening class and so on. This course sequence will include the following four years (assuming the user is a student who just started his </w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
        '<w:r><w:t xml:space="preserve">or her undergraduate program). Once the course sequence is created, a student will be able to access it to view the generated schedule by the </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>TimeTurner</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>.</w:t></w:r>' +
        '</w:p>' +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$tpr.InsertXML($xml3)
